$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or "208.65" are not coerced to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "79.572.16"
$ws.Range("E2").Value = "  +4.20%  "
$ws.Range("D3").Value = "3.198.05"
$ws.Range("E3").Value = "  +5.25%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "208.65"
$ws.Range("E5").Value = "  +5.35%  "
$ws.Range("D6").Value = "633.83"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.245"
$ws.Range("E8").Value = "  +19.52%  "
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +9.59%  "
$ws.Range("D10").Value = "3.195.64"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").Value = "0.604"
$ws.Range("E11").Value = "  +38.65%  "
$ws.Range("D12").Value = "0.0000261"
$ws.Range("E12").Value = "  +35.87%  "
$ws.Range("D13").Value = "0.165"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").Value = "5.39"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("D15").Value = "3.786.05"
$ws.Range("E15").Value = "  +5.23%  "
$ws.Range("D16").Value = "31.99"
$ws.Range("E16").Value = "  +11.24%  "
$ws.Range("D17").Value = "79.407.91"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("D18").Value = "3.192.27"
$ws.Range("E18").Value = "  +5.22%  "
$ws.Range("D19").Value = "14.50"
$ws.Range("E19").Value = "  +7.41%  "
$ws.Range("D20").Value = "9.35"
$ws.Range("E20").Value = "  +4.28%  "
$ws.Range("D21").Value = "2.99"
$ws.Range("E21").Value = "  +27.08%  "
$ws.Range("D22").Value = "441.13"
$ws.Range("E22").Value = "  +16.10%  "
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").Value = "  +20.29%  "
$ws.Range("D24").Value = "4.82"
$ws.Range("E24").Value = "  +11.26%  "
$ws.Range("D25").Value = "3.364.08"
$ws.Range("E25").Value = "  +5.50%  "
$ws.Range("D26").Value = "77.17"
$ws.Range("E26").Value = "  +6.41%  "
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  +11.61%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "0.0000124"
$ws.Range("E29").Value = "  +15.62%  "
$ws.Range("D30").Value = "9.12"
$ws.Range("E30").Value = "  +10.50%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "1.52"
$ws.Range("E32").Value = "  +9.30%  "
$ws.Range("D33").Value = "548.66"
$ws.Range("E33").Value = "  +11.17%  "
$ws.Range("D34").Value = "0.155"
$ws.Range("E34").Value = "  +32.35%  "
$ws.Range("D35").Value = "2.03"
$ws.Range("E35").Value = "  +5.89%  "
$ws.Range("D36").Value = "23.02"
$ws.Range("E36").Value = "  +11.98%  "
$ws.Range("D37").Value = "0.123"
$ws.Range("E37").Value = "  +17.92%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "0.412"
$ws.Range("E39").Value = "  +8.88%  "
$ws.Range("D40").Value = "163.33"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "20.03"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "191.70"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "5.56"
$ws.Range("E44").Value = "  +9.92%  "
$ws.Range("D45").Value = "1.83"
$ws.Range("E45").Value = "  +11.78%  "
$ws.Range("D46").Value = "0.800"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").Value = "2.67"
$ws.Range("E47").Value = "  +11.06%  "
$ws.Range("D48").Value = "1.34"
$ws.Range("E48").Value = "  +6.35%  "
$ws.Range("D49").Value = "43.09"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").Value = "0.641"
$ws.Range("E50").Value = "  +7.08%  "
$ws.Range("D51").Value = "25.54"
$ws.Range("E51").Value = "  +15.48%  "
# Restore the original (default) cell style now that the text values are safely stored
$ws.Range("D2:D51").Style = "Normal"
